$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price/Volume(1h) cells are stored as plain text (not numbers or
# percentages) in the source workbook. Assigning a bare numeric-looking
# string lets Excel auto-convert it to a real number/percentage, so we
# prefix with an apostrophe to force text entry, then reset the cell
# style back to Normal so no stray quote-prefix / text-format style ends
# up attached to the cell (keeping formatting identical to the original).
$ws.Range("D2").Value = "'259.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.53%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'-0.98%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.689"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.23%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06027"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.36%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.678"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.50%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8588"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.10%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9303"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.29%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-1.04%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.04833"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'21.34%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07060"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.36%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03148"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.90%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09130"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.29%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001527"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.17%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006046"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.32%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005997"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-3.12%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.462"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.24%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'-1.25%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E20").Value = "'0.45%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'0.40%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.113"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'5.51%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04235"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.34%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001215"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.49%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004043"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "'-0.08%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-21.37%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03845"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.06%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1114"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.07%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.003938"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.46%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002292"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-4.86%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'29.84%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005096"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-6.54%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.06%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.05449"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-9.13%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'0.95%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.06%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("E50").Style = "Normal"
